$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns for the rows being edited so that
# numeric-looking strings (e.g. "0.571", "1.00") are kept as literal text, matching
# the original inline-string cell contents, instead of being auto-converted to numbers.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = '63.830.54'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '3.419.66'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("D5").Value = '570.36'
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("D6").Value = '157.20'
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.421.31'
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '0.571'
$ws.Range("E9").Value = '  -6.42%  '
$ws.Range("D10").Value = '7.25'
$ws.Range("E10").Value = '  +0.86%  '
$ws.Range("D11").Value = '0.120'
$ws.Range("E11").Value = '  -1.77%  '
$ws.Range("E12").Value = '  -2.84%  '
$ws.Range("D13").Value = '4.007.40'
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("D15").Value = '27.19'
$ws.Range("E15").Value = '  -2.40%  '
$ws.Range("E16").Value = '  -6.99%  '
$ws.Range("D17").Value = '63.894.44'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").Value = '3.429.16'
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("E19").Value = '  -3.54%  '
$ws.Range("D20").Value = '13.65'
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("D21").Value = '383.16'
$ws.Range("E21").Value = '  +2.43%  '
$ws.Range("D22").Value = '7.77'
$ws.Range("E22").Value = '  -2.25%  '
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("D24").Value = '71.19'
$ws.Range("E24").Value = '  -0.96%  '
$ws.Range("D25").Value = '0.520'
$ws.Range("E25").Value = '  -5.13%  '
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("D27").Value = '9.68'
$ws.Range("E27").Value = '  -3.34%  '
$ws.Range("E28").Value = '  +1.00%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = '6.13'
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("E31").Value = '  -5.12%  '
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("D33").Value = '22.99'
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("D35").Value = '7.01'
$ws.Range("E35").Value = '  -2.88%  '
$ws.Range("E36").Value = '  -4.21%  '
$ws.Range("D37").Value = '161.20'
$ws.Range("E37").Value = '  +0.60%  '
$ws.Range("D38").Value = '0.833'
$ws.Range("E38").Value = '  +8.77%  '
$ws.Range("D39").Value = '1.84'
$ws.Range("E39").Value = '  -2.41%  '
$ws.Range("E40").Value = '  -1.93%  '
$ws.Range("D41").Value = '2.812.07'
$ws.Range("E41").Value = '  -1.08%  '
$ws.Range("D42").Value = '0.0726'
$ws.Range("E42").Value = '  -3.88%  '
$ws.Range("D43").Value = '42.82'
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("D44").Value = '4.43'
$ws.Range("E44").Value = '  -3.86%  '
$ws.Range("D45").Value = '6.40'
$ws.Range("E45").Value = '  -5.01%  '
$ws.Range("D46").Value = '25.85'
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").Value = '2.33'
$ws.Range("E48").Value = '  +9.54%  '
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").Value = '328.07'
$ws.Range("E49").Value = '  +4.50%  '
$ws.Range("E50").Value = '  -3.53%  '
$ws.Range("E51").Value = '  -5.02%  '

# Restore the original (default) style so no residual formatting diff is left behind.
$editRange.Style = "Normal"
